$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.122774
$ws.Range("H2").Value = 33.368322
$ws.Range("I2").Value = 0.2449652610853511
$ws.Range("J2").Value = 0.2449652610853511
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 1707.668562240749
$ws.Range("R2").Value = 15369.01706016674
$ws.Range("S2").Value = 0.07770805083333913
$ws.Range("T2").Value = 0.07770805083333915

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.122774
$ws.Range("H3").Value = 33.368322
$ws.Range("I3").Value = 0.2449652610853511
$ws.Range("J3").Value = 0.2449652610853511
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 1877.520944028531
$ws.Range("R3").Value = 16897.68849625678
$ws.Range("S3").Value = 0.08543724243994077
$ws.Range("T3").Value = 0.08543724243994077

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.122774
$ws.Range("H4").Value = 33.368322
$ws.Range("I4").Value = 0.2449652610853511
$ws.Range("J4").Value = 0.2449652610853511
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 757.3532780235932
$ws.Range("R4").Value = 6816.17950221234
$ws.Range("S4").Value = 0.0344636238722045
$ws.Range("T4").Value = 0.03446362387220451

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.122774
$ws.Range("H5").Value = 33.368322
$ws.Range("I5").Value = 0.2449652610853511
$ws.Range("J5").Value = 0.2449652610853511
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 1040.676466614902
$ws.Range("R5").Value = 9366.088199534117
$ws.Range("S5").Value = 0.04735634393986669
$ws.Range("T5").Value = 0.04735634393986669

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.220714
$ws.Range("H6").Value = 54.662142
$ws.Range("I6").Value = 0.4012885600454987
$ws.Range("J6").Value = 0.4012885600454988
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 2797.408315531709
$ws.Range("R6").Value = 25176.67483978538
$ws.Range("S6").Value = 0.1272970366683468
$ws.Range("T6").Value = 0.1272970366683468

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.220714
$ws.Range("H7").Value = 54.662142
$ws.Range("I7").Value = 0.4012885600454987
$ws.Range("J7").Value = 0.4012885600454988
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 3075.651105574371
$ws.Range("R7").Value = 27680.85995016934
$ws.Range("S7").Value = 0.1399585714361204
$ws.Range("T7").Value = 0.1399585714361204

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.220714
$ws.Range("H8").Value = 54.662142
$ws.Range("I8").Value = 0.4012885600454987
$ws.Range("J8").Value = 0.4012885600454988
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 1240.654307624193
$ws.Range("R8").Value = 11165.88876861774
$ws.Range("S8").Value = 0.05645640502800927
$ws.Range("T8").Value = 0.05645640502800928

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.220714
$ws.Range("H9").Value = 54.662142
$ws.Range("I9").Value = 0.4012885600454987
$ws.Range("J9").Value = 0.4012885600454988
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 1704.778705808522
$ws.Range("R9").Value = 15343.0083522767
$ws.Range("S9").Value = 0.07757654691302227
$ws.Range("T9").Value = 0.07757654691302227

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.1189986666666667
$ws.Range("H10").Value = 0.356996
$ws.Range("I10").Value = 0.002620797603979787
$ws.Range("J10").Value = 0.002620797603979787
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 18.26974835731022
$ws.Range("R10").Value = 164.427735215792
$ws.Range("S10").Value = 0.0008313712423207477
$ws.Range("T10").Value = 0.0008313712423207479

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.1189986666666667
$ws.Range("H11").Value = 0.356996
$ws.Range("I11").Value = 0.002620797603979787
$ws.Range("J11").Value = 0.002620797603979787
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 20.08693955106311
$ws.Range("R11").Value = 180.782455959568
$ws.Range("S11").Value = 0.0009140631585276927
$ws.Range("T11").Value = 0.0009140631585276927

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.1189986666666667
$ws.Range("H12").Value = 0.356996
$ws.Range("I12").Value = 0.002620797603979787
$ws.Range("J12").Value = 0.002620797603979787
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 8.102657689568888
$ws.Range("R12").Value = 72.92391920612
$ws.Range("S12").Value = 0.000368714251435284
$ws.Range("T12").Value = 0.000368714251435284

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.1189986666666667
$ws.Range("H13").Value = 0.356996
$ws.Range("I13").Value = 0.002620797603979787
$ws.Range("J13").Value = 0.002620797603979787
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 11.13383333676933
$ws.Range("R13").Value = 100.204500030924
$ws.Range("S13").Value = 0.0005066489516960622
$ws.Range("T13").Value = 0.0005066489516960622

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 15.943029
$ws.Range("H14").Value = 47.829087
$ws.Range("I14").Value = 0.3511253812651704
$ws.Range("J14").Value = 0.3511253812651704
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 2447.717575686836
$ws.Range("R14").Value = 22029.45818118152
$ws.Range("S14").Value = 0.1113842381378422
$ws.Range("T14").Value = 0.1113842381378423

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 15.943029
$ws.Range("H15").Value = 47.829087
$ws.Range("I15").Value = 0.3511253812651704
$ws.Range("J15").Value = 0.3511253812651704
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 2691.178554806045
$ws.Range("R15").Value = 24220.6069932544
$ws.Range("S15").Value = 0.1224630145231762
$ws.Range("T15").Value = 0.1224630145231762

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 15.943029
$ws.Range("H16").Value = 47.829087
$ws.Range("I16").Value = 0.3511253812651704
$ws.Range("J16").Value = 0.3511253812651704
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 1085.56599952271
$ws.Range("R16").Value = 9770.09399570439
$ws.Range("S16").Value = 0.04939905772064133
$ws.Range("T16").Value = 0.04939905772064134

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 15.943029
$ws.Range("H17").Value = 47.829087
$ws.Range("I17").Value = 0.3511253812651704
$ws.Range("J17").Value = 0.3511253812651704
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 1491.672408956517
$ws.Range("R17").Value = 13425.05168060865
$ws.Range("S17").Value = 0.06787907088351063
$ws.Range("T17").Value = 0.06787907088351063
